$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 321; this shifts old rows 321-394 down to 322-395,
# matching the target diff (dimension grows from R394 to R395).
$ws.Rows.Item(321).Insert()

# Populate the newly inserted row 321 with the new record.
$ws.Range("A321").Value = 4
$ws.Range("B321").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C321").Value = "Los Lagos"
$ws.Range("D321").Value = 44711
$ws.Range("E321").Value = 10
$ws.Range("F321").Value = 100114001
$ws.Range("G321").Value = "Papa"
$ws.Range("H321").Value = "Patagonia"
$ws.Range("I321").Value = "1a (guarda)"
$ws.Range("J321").Value = 250
$ws.Range("K321").Value = 7000
$ws.Range("L321").Value = 7500
$ws.Range("M321").Value = 7200
$ws.Range("N321").Value = "`$/saco 25 kilos"
$ws.Range("O321").Value = "Provincia de Llanquihue"
$ws.Range("P321").Value = 288
$ws.Range("Q321").Value = 25
$ws.Range("R321").Value = "Hortaliza"
